$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'30.116.53"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +0.02%  '

$ws.Range("D3").Formula = "'1.884.43"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.21%  '

$ws.Range("D4").Formula = "'0.9975"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.10%  '

$ws.Range("D5").Formula = "'244.34"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.75%  '

$ws.Range("D6").Formula = "'0.9981"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.08%  '

$ws.Range("D7").Formula = "'0.4976"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.02%  '

$ws.Range("D8").Formula = "'44.51"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.84%  '

$ws.Range("D9").Formula = "'0.2920"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.30%  '

$ws.Range("D10").Formula = "'0.06631"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.75%  '

$ws.Range("D11").Formula = "'1.879.09"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.06%  '

$ws.Range("D12").Formula = "'16.88"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.54%  '

$ws.Range("D13").Formula = "'0.07207"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.28%  '

$ws.Range("D14").Formula = "'0.6661"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.29%  '

$ws.Range("D15").Formula = "'85.73"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.58%  '

$ws.Range("D16").Formula = "'4.837"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +0.70%  '

$ws.Range("D17").Formula = "'30.092.29"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.06%  '

$ws.Range("D18").Formula = "'0.000007840"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +4.10%  '

$ws.Range("D19").Formula = "'0.9976"
$ws.Range("D19").Style = "Normal"

$ws.Range("E20").Value = '  -0.48%  '

$ws.Range("D21").Formula = "'2.122.09"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.38%  '

$ws.Range("D22").Formula = "'0.9980"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.01%  '

$ws.Range("D23").Formula = "'4.766"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.53%  '

$ws.Range("B24").Value = 'Cosmos'
$ws.Range("C24").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D24").Formula = "'9.175"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.48%  '

$ws.Range("B25").Value = 'Chainlink'
$ws.Range("C25").Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range("D25").Formula = "'5.609"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.70%  '

$ws.Range("D26").Formula = "'151.17"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.33%  '

$ws.Range("D27").Formula = "'135.34"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.15%  '

$ws.Range("E28").Value = '  +0.46%  '

$ws.Range("E29").Value = '  -2.35%  '

$ws.Range("D30").Formula = "'1.380"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.45%  '

$ws.Range("D31").Formula = "'4.166"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.86%  '

$ws.Range("D32").Formula = "'0.08676"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +0.92%  '

$ws.Range("D33").Formula = "'3.952"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.74%  '

$ws.Range("D34").Formula = "'0.05005"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -1.07%  '

$ws.Range("B35").Value = 'ImmutableX'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D35").Formula = "'0.7062"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +3.41%  '

$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").Formula = "'1.107"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -2.75%  '

$ws.Range("D37").Formula = "'2.658"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.79%  '

$ws.Range("D38").Formula = "'2.703"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -1.37%  '

$ws.Range("D39").Formula = "'2.200"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -5.08%  '

$ws.Range("D40").Formula = "'0.9347"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -2.93%  '

$ws.Range("D41").Formula = "'0.01649"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.56%  '

$ws.Range("D42").Formula = "'5.964"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -1.87%  '

$ws.Range("D43").Formula = "'0.9990"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.06%  '

$ws.Range("D44").Formula = "'0.4205"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.01%  '

$ws.Range("D45").Formula = "'101.81"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -1.24%  '

$ws.Range("D46").Formula = "'7.515"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.91%  '

$ws.Range("D47").Formula = "'0.1261"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +0.31%  '

$ws.Range("D48").Formula = "'0.05715"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.71%  '

$ws.Range("D49").Formula = "'32.45"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.04%  '

$ws.Range("D50").Formula = "'8.283"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.06%  '

$ws.Range("D51").Formula = "'1.350"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.88%  '
